# Updated symbol list on Wed Jan 18 17:11:20 UTC 2023 with GitHub Actions
#
# The crypto price/volume snapshot was refreshed: prices, 1h volume %
# changes and the "Hora" (hour) column were updated for the new
# 17:00 scrape, and a newly-tracked coin (GateToken) pushed into rank
# #8, shifting the coins that were previously ranked #8-#17 down by
# one row each (their price/volume values were also refreshed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. Values are written as Text so
# that numeric-looking strings (prices, "17", "2,116.77%", etc.) are
# preserved exactly as text, matching the workbook's original
# inline-string cell contents instead of being coerced into numbers.
$updates = [ordered]@{
    "D2" = "294.90"
    "E2" = "-2.50%"
    "G2" = "17"
    "D3" = "31.32"
    "E3" = "-1.68%"
    "G3" = "17"
    "D4" = "4.983"
    "E4" = "-0.39%"
    "G4" = "17"
    "D5" = "0.07362"
    "E5" = "-5.75%"
    "G5" = "17"
    "D6" = "1.838"
    "E6" = "-13.02%"
    "G6" = "17"
    "D7" = "7.620"
    "E7" = "-2.18%"
    "G7" = "17"
    "B8" = "GateToken"
    "C8" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "D8" = "3.742"
    "E8" = "-1.04%"
    "G8" = "17"
    "B9" = "MXToken"
    "C9" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "D9" = "0.9144"
    "E9" = "-0.85%"
    "G9" = "17"
    "B10" = "WazirX"
    "C10" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D10" = "0.1646"
    "E10" = "-5.92%"
    "G10" = "17"
    "B11" = "LiechtensteinCryptoassetsExchange"
    "C11" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "D11" = "0.07584"
    "E11" = "-2.97%"
    "G11" = "17"
    "B12" = "MandalaExchangeToken"
    "C12" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "D12" = "0.08174"
    "E12" = "-7.32%"
    "G12" = "17"
    "B13" = "BitrueCoin"
    "C13" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "D13" = "0.02998"
    "E13" = "-2.90%"
    "G13" = "17"
    "B14" = "BitMartToken"
    "C14" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D14" = "0.1006"
    "E14" = "0.72%"
    "G14" = "17"
    "B15" = "BitForexToken"
    "C15" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D15" = "0.001497"
    "E15" = "-0.87%"
    "G15" = "17"
    "B16" = "TigerCash"
    "C16" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "D16" = "0.005696"
    "E16" = "-2.25%"
    "G16" = "17"
    "B17" = "UpBots"
    "C17" = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
    "D17" = "0.007498"
    "E17" = "2,116.77%"
    "G17" = "17"
    "B18" = "LEO"
    "C18" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D18" = "3.463"
    "E18" = "-0.03%"
    "G18" = "17"
    "D19" = "2.127"
    "E19" = "-6.25%"
    "G19" = "17"
    "D20" = "0.3262"
    "E20" = "-0.92%"
    "G20" = "17"
    "D21" = "0.1291"
    "E21" = "-2.78%"
    "G21" = "17"
    "D22" = "4.365"
    "E22" = "4.86%"
    "G22" = "17"
    "D23" = "0.1973"
    "E23" = "9.22%"
    "G23" = "17"
    "D24" = "0.04492"
    "E24" = "-2.14%"
    "G24" = "17"
    "D25" = "0.001233"
    "E25" = "-0.50%"
    "G25" = "17"
    "D26" = "0.004044"
    "E26" = "-9.62%"
    "G26" = "17"
    "D27" = "0.0001259"
    "E27" = "0.60%"
    "G27" = "17"
    "G28" = "17"
    "G29" = "17"
    "G30" = "17"
    "G31" = "17"
    "G32" = "17"
    "G33" = "17"
    "G34" = "17"
    "G35" = "17"
    "G36" = "17"
    "G37" = "17"
    "G38" = "17"
    "D39" = "0.01614"
    "E39" = "-7.97%"
    "G39" = "17"
    "D40" = "0.04415"
    "E40" = "-7.30%"
    "G40" = "17"
    "D41" = "0.007432"
    "E41" = "4.39%"
    "G41" = "17"
    "D42" = "0.1328"
    "E42" = "-2.05%"
    "G42" = "17"
    "D43" = "0.002070"
    "E43" = "-3.35%"
    "G43" = "17"
    "D44" = "0.01104"
    "E44" = "2.40%"
    "G44" = "17"
    "D45" = "0.00006044"
    "E45" = "1.19%"
    "G45" = "17"
    "D46" = "0.00000000755"
    "E46" = "0.69%"
    "G46" = "17"
    "D47" = "1.777"
    "E47" = "51.18%"
    "G47" = "17"
    "D48" = "0.002907"
    "E48" = "-18.08%"
    "G48" = "17"
    "D49" = "0.00002113"
    "E49" = "0.69%"
    "G49" = "17"
    "D50" = "0.0002013"
    "E50" = "0.69%"
    "G50" = "17"
    "G51" = "17"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}
